$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. "Förändrad" (changed) date column C: all existing rows 2-501 move
#    from 45192 (2023-09-12) to 45202 (2023-09-22).
$ws.Range("C2:C501").Value = 45202

# 2. Row 501 gains an explicit (default) row height marker.
$ws.Range("A501").EntireRow.RowHeight = 15

# 3. Two new cutting-notice rows appended at the bottom of the table.

# Row 502
$ws.Range("A502").Value = "A 45865-2023"
$ws.Range("B502").Value = 45195
$ws.Range("B502").NumberFormat = "YYYY-MM-DD"
$ws.Range("C502").Value = 45202
$ws.Range("C502").NumberFormat = "YYYY-MM-DD"
$ws.Range("D502").Value = "JÖNKÖPINGS LÄN"
$ws.Range("E502").Value = "NÄSSJÖ"
$ws.Range("G502").Value = 2.6
$ws.Range("H502").Value = 0
$ws.Range("I502").Value = 0
$ws.Range("J502").Value = 0
$ws.Range("K502").Value = 0
$ws.Range("L502").Value = 0
$ws.Range("M502").Value = 0
$ws.Range("N502").Value = 0
$ws.Range("O502").Value = 0
$ws.Range("P502").Value = 0
$ws.Range("Q502").Value = 0
$ws.Range("R502").WrapText = $true
$ws.Range("A502").EntireRow.RowHeight = 15

# Row 503
$ws.Range("A503").Value = "A 45868-2023"
$ws.Range("B503").Value = 45195
$ws.Range("B503").NumberFormat = "YYYY-MM-DD"
$ws.Range("C503").Value = 45202
$ws.Range("C503").NumberFormat = "YYYY-MM-DD"
$ws.Range("D503").Value = "JÖNKÖPINGS LÄN"
$ws.Range("E503").Value = "NÄSSJÖ"
$ws.Range("G503").Value = 1.7
$ws.Range("H503").Value = 0
$ws.Range("I503").Value = 0
$ws.Range("J503").Value = 0
$ws.Range("K503").Value = 0
$ws.Range("L503").Value = 0
$ws.Range("M503").Value = 0
$ws.Range("N503").Value = 0
$ws.Range("O503").Value = 0
$ws.Range("P503").Value = 0
$ws.Range("Q503").Value = 0
$ws.Range("R503").WrapText = $true
